$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings
# (e.g. "0.9964", "317.13") are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.344.85"
$ws.Range("E2").Value = "  +4.88%  "
$ws.Range("D3").Value = "1.811.65"
$ws.Range("E3").Value = "  +4.88%  "
$ws.Range("D4").Value = "0.9964"
$ws.Range("D5").Value = "317.13"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "0.9959"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "0.5664"
$ws.Range("E7").Value = "  +17.81%  "
$ws.Range("D8").Value = "0.3843"
$ws.Range("E8").Value = "  +10.40%  "
$ws.Range("D9").Value = "0.07642"
$ws.Range("E9").Value = "  +5.59%  "
$ws.Range("D10").Value = "43.26"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "1.137"
$ws.Range("E11").Value = "  +8.28%  "
$ws.Range("D12").Value = "21.37"
$ws.Range("E12").Value = "  +7.27%  "
$ws.Range("D13").Value = "0.9944"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "6.238"
$ws.Range("E14").Value = "  +6.21%  "
$ws.Range("D15").Value = "1.799.84"
$ws.Range("E15").Value = "  +4.61%  "
$ws.Range("D16").Value = "7.226"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").Value = "92.26"
$ws.Range("E17").Value = "  +6.24%  "
$ws.Range("D18").Value = "0.00001080"
$ws.Range("D19").Value = "0.06507"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "0.9957"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").Value = "5.992"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").Value = "28.367.41"
$ws.Range("E23").Value = "  +4.76%  "
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").Value = "2.101"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "20.80"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").Value = "156.62"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "2.374"
$ws.Range("E28").Value = "  +14.89%  "
$ws.Range("D29").Value = "2.009.48"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("D30").Value = "123.55"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "1.152"
$ws.Range("E31").Value = "  +10.21%  "
$ws.Range("D32").Value = "0.1043"
$ws.Range("E32").Value = "  +12.01%  "
$ws.Range("D33").Value = "5.762"
$ws.Range("E33").Value = "  +7.19%  "
$ws.Range("D34").Value = "3.612"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "0.02314"
$ws.Range("E35").Value = "  +5.93%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "11.67"
$ws.Range("E36").Value = "  +6.32%  "
$ws.Range("D37").Value = "8.688"
$ws.Range("E37").Value = "  +15.30%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2127"
$ws.Range("E38").Value = "  +6.75%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.6432"
$ws.Range("E39").Value = "  +7.62%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.033"
$ws.Range("E40").Value = "  +5.64%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.06062"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "0.9956"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.153"
$ws.Range("E43").Value = "  +5.30%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.390"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  +5.56%  "
$ws.Range("D46").Value = "0.5980"
$ws.Range("E46").Value = "  +6.40%  "
$ws.Range("D47").Value = "3.692"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "122.67"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "1.938"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "1.141"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").Value = "0.06832"
$ws.Range("E51").Value = "  +3.11%  "
